$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NIG(0.8359802999628896, 0.580617726797404, 2.0397646119585824, 3.44681798555417)"
$ws.Range("C2").Value = "JSU(-1.5585472062144725, 1.6681034567163944, -0.21802181981375768, 9.249868129668949)"
$ws.Range("D2").Value = "JSU(-1.0051392706149738, 1.1324737740393926, 0.4129442553597543, 2.7444956726146312)"
$ws.Range("E2").Value = "NIG(0.908248506523404, 0.6069449054983979, 5.1692827924436955, 5.3454507061725)"
